# ---------------------------------------------------------------------------
# Edit script for VT_fee_waiver_next_steps.docx
# Implements the changes described by the diff:
#  1. Insert a "{%p if trial_court_division != "Unknown" %}" paragraph and an
#     empty paragraph before the "Details for the court you selected:"
#     paragraph; give that paragraph w:after="240" spacing.
#  2. Remove the stray empty paragraph that used to follow it.
#  3. Split the paragraph ending in "{{ phone }}" so a new
#     "{% p endif %}" paragraph (carrying the trailing <w:br/> run) follows it.
#  4. Move <w:lastRenderedPageBreak/> from before "to pay " to the start of
#     the "If the waiver is " paragraph, and merge the two runs back
#     together.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument
$wdNS = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function New-EmptyParaXml {
    return "<w:p $wdNS></w:p>"
}

# ---------------------------------------------------------------------------
# Hunk 1: insert two paragraphs before "Details for the court you selected:"
# ---------------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("e-filing system", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$eFileListPara = $rng.Paragraphs(1)

$insertPoint = $d.Range($eFileListPara.Range.End, $eFileListPara.Range.End)

$hunk1Xml = @"
<w:p $wdNS>
  <w:r><w:t>{%</w:t></w:r>
  <w:r><w:t xml:space="preserve">p </w:t></w:r>
  <w:r><w:t xml:space="preserve">if </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>trial_court_division</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t>!</w:t></w:r>
  <w:r><w:t>= "</w:t></w:r>
  <w:r><w:t>Unknown</w:t></w:r>
  <w:r><w:t>"</w:t></w:r>
  <w:r><w:t xml:space="preserve"> %}</w:t></w:r>
</w:p>
<w:p $wdNS></w:p>
$(New-EmptyParaXml)
"@

$insertPoint.InsertXML($hunk1Xml)

# Remove the trailing stray empty paragraph produced by InsertXML's merge
# semantics (it sits directly before "Details for the court you selected:").
$rng = $d.Content
$rng.Find.Execute("Details for the court you selected", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$detailsPara = $rng.Paragraphs(1)
$strayPara = $detailsPara.Previous()
$strayPara.Range.Delete()

# Set spacing for the two newly-added paragraphs (w:after=0, single line).
$rng = $d.Content
$rng.Find.Execute("e-filing system", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$eFileListPara = $rng.Paragraphs(1)
$ifPara = $eFileListPara.Next()
$ifPara.Format.SpaceAfter = 0
$ifPara.Format.LineSpacingRule = 0
$blankPara = $ifPara.Next()
$blankPara.Format.SpaceAfter = 0
$blankPara.Format.LineSpacingRule = 0

# Give the "Details..." paragraph its new spacing (w:after=240).
$rng = $d.Content
$rng.Find.Execute("Details for the court you selected", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$detailsPara = $rng.Paragraphs(1)
$detailsPara.Format.SpaceAfter = 12

# Remove the now-obsolete empty paragraph that used to follow "Details...:"
$oldEmptyPara = $detailsPara.Next()
$oldEmptyPara.Range.Delete()

# ---------------------------------------------------------------------------
# Hunk 2: split the "...{{ phone }}<br/>" paragraph so the trailing <w:br/>
# run moves into its own new paragraph, preceded by "{% p endif %}".
# ---------------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("{{ phone }}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.InsertParagraphAfter()

$rng = $d.Content
$rng.Find.Execute("{{ phone }}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$phonePara = $rng.Paragraphs(1)
$brPara = $phonePara.Next()

$insertPt = $d.Range($brPara.Range.Start, $brPara.Range.Start)
$endifFragment = "<w:p $wdNS>" +
    "<w:r><w:t>{%</w:t></w:r>" +
    "<w:r><w:t>p</w:t></w:r>" +
    '<w:r><w:t xml:space="preserve"> endif %}</w:t></w:r>' +
    "</w:p>"
$insertPt.InsertXML($endifFragment)
